$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) column before column N ---
# This shifts the old N/O/P (Late / heading / Outstanding) columns to O/P/Q.
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns.Item(14).Insert() | Out-Null

# The newly inserted column N keeps a specific custom width (no bestFit),
# matching the other data columns on this sheet.
$ws.Columns.Item(14).ColumnWidth = 10.166666666666666

# Update the sheet's current selection to S6.
$ws.Range("S6").Select() | Out-Null

# Make "Repayment schedule" the active/selected tab of the workbook
# (was previously "Transactions").
$ws.Select() | Out-Null
